$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Step 1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2 content edits ("set 3 w/o ES") ---
# Widen column C to fit the longer note text (closest value the pixel-quantized
# ColumnWidth property can reach to the target 48.7109375 stored width)
$ws2.Columns.Item(3).ColumnWidth = 47.86

# Row 3 (Set 3): accuracy cell becomes an annotated note instead of a bare number
$ws2.Range("C3").Value = "0.9469517023 (need to redo - without early stopping)"

# Row 5 (Set 4): add the same "redo" note
$ws2.Range("C5").Value = "need to redo - without early stopping"

# Row 6 (Set 5): fill in the test accuracy + epoch count that were missing
$ws2.Range("C6").Value = 0.90981789390340395
$ws2.Range("H6").Value = 45

# --- Selection / active sheet bookkeeping ---
# Previously Sheet2 was the active tab with C4 selected; now the user has
# moved back to "Step 1" (selecting C5) leaving Sheet2's selection at B6.
$ws2.Activate()
$ws2.Range("B6").Select()

$ws1.Activate()
$ws1.Range("C5").Select()

Write-Output "done"
